$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '24.464.25'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -1.67%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.656.42'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -3.17%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '306.74'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -1.23%  '
$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9989'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +0.20%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3620'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -3.42%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '47.10'
$c.Style = "Normal"
$ws.Range("E8").Value = '  -4.95%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.3241'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -5.86%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.123'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -7.07%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07031'
$c.Style = "Normal"
$ws.Range("E11").Value = '  -6.89%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.9995'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.11%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '5.903'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -6.60%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '19.43'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -8.17%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '1.657.80'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -3.12%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '6.577'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -6.58%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.00001047'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -7.92%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.06581'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.54%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.9988'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.19%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '77.67'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -8.40%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '5.917'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -7.48%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '15.65'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -9.73%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '12.48'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -5.25%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '24.464.35'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.62%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.468'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.74%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.342'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -16.17%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '146.31'
$c.Style = "Normal"
$ws.Range("E27").Value = '  -3.49%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '18.58'
$c.Style = "Normal"
$ws.Range("E28").Value = '  -8.97%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.842.65'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -3.03%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '124.49'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -5.79%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.182'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -5.10%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '4.033'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -5.99%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '5.712'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -17.94%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.08443'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -4.36%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.665'
$c.Style = "Normal"
$ws.Range("E35").Value = '  -10.42%  '
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '12.32'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -11.15%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.223'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -7.11%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.06031'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -9.67%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.02208'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -8.56%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.2059'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -8.53%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.199'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -6.21%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '8.168'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -11.07%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.9987'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +0.18%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.5906'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -8.93%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.730'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -2.49%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '12.67'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -9.06%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.5621'
$c.Style = "Normal"
$ws.Range("E47").Value = '  -9.21%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '122.30'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -6.31%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.939'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -9.55%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.06925'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -5.41%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '74.31'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -7.19%  '